$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a "Name" column (B) and a "Date" column (F) to the jobs table.
# Existing Client/Type data shifts from B/C into C/D; a new "Name" value
# (previously hardcoded in column D) now lives in B; File (E) is
# untouched; a new "Date" column is filled in at F.

$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Client"
$ws.Range("D1").Value = "Type"
$ws.Range("F1").Value = "Date"

$dateValue = "Feb 29, 2004 (00:00:00 EST)"

$rows = @(
    @{ Row = 2;  Name = "Supplier";     Client = "DRX"; Type = "I" },
    @{ Row = 3;  Name = "Plant";        Client = "DRX"; Type = "I" },
    @{ Row = 4;  Name = "Solicitation"; Client = "DRX"; Type = "I" },
    @{ Row = 5;  Name = "BOM";          Client = "DRX"; Type = "I" },
    @{ Row = 6;  Name = "RequestFile";  Client = "DRX"; Type = "E" },
    @{ Row = 7;  Name = "BOM";          Client = "DRX"; Type = "E" },
    @{ Row = 8;  Name = "Supplier";     Client = "GYU"; Type = "I" },
    @{ Row = 9;  Name = "Plant";        Client = "GYU"; Type = "I" },
    @{ Row = 10; Name = "Solicitation"; Client = "GYU"; Type = "I" },
    @{ Row = 11; Name = "BOM";          Client = "GYU"; Type = "I" },
    @{ Row = 12; Name = "RequestFile";  Client = "GYU"; Type = "E" },
    @{ Row = 13; Name = "BOM";          Client = "GYU"; Type = "E" }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("B$n").Value = $r.Name
    $ws.Range("C$n").Value = $r.Client
    $ws.Range("D$n").Value = $r.Type
    $ws.Range("F$n").Value = $dateValue
}

$ws.Range("H6").Select()
